$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.795.92'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -4.17%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.159.11'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -4.40%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '525.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.38'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.11%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.156.74'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.452'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -6.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.30'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.63%  '
$ws.Range("E11").Value = '  -6.80%  '
$ws.Range("E12").Value = '  -3.89%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.702.64'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.37%  '
$ws.Range("E14").Value = '  -1.88%  '
$ws.Range("E15").Value = '  -4.95%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.159.73'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.74%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '57.790.81'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.24%  '
$ws.Range("E18").Value = '  -7.72%  '
$ws.Range("E19").Value = '  -5.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.06'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -8.91%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.06'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '346.54'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -7.05%  '
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.64'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.512'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.292.11'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.12%  '
$ws.Range("E27").Value = '  -8.49%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.166'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.60%  '
$ws.Range("E29").Value = '  +0.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.86'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.75%  '
$ws.Range("E31").Value = '  -0.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.89'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.81%  '
$ws.Range("E33").Value = '  -9.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '21.74'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.60%  '
$ws.Range("E35").Value = '  -4.82%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.90'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '160.12'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.25'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.40'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.68%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '25.79'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0697'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.186.02'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.60%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.63'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.700'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.08'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.95'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.05%  '
$ws.Range("E47").Value = '  -0.34%  '
$ws.Range("E48").Value = '  -8.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.269.70'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.18'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.54'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.03%  '
